$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4849.6294
$ws.Range("I40").Value = 1746.6666
$ws.Range("J40").Value = 8728.333000000001
$ws.Range("K40").Value = 1746.6666
$ws.Range("L40").Value = 8728.333000000001
$ws.Range("M40").Value = -1571.6666
$ws.Range("N40").Value = -9078.333000000001

$ws.Range("H43").Value = 2874.25
$ws.Range("I43").Value = 2998
$ws.Range("J43").Value = 2833
$ws.Range("K43").Value = 2998
$ws.Range("L43").Value = 2833
$ws.Range("M43").Value = -2929
$ws.Range("N43").Value = -2971

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H62").Value = 4100
$ws.Range("I62").Value = 3800
$ws.Range("K62").Value = 3800
$ws.Range("M62").Value = -3176

$ws.Range("H65").Value = 4100
$ws.Range("I65").Value = 3800
$ws.Range("K65").Value = 19000
$ws.Range("M65").Value = -15880

$ws.Range("H96").Value = 419.45456
$ws.Range("I96").Value = 271
$ws.Range("J96").Value = 1087.5
$ws.Range("K96").Value = 813
$ws.Range("L96").Value = 3262.5
$ws.Range("M96").Value = 560
$ws.Range("N96").Value = -6008.5

$ws.Range("H106").Value = 5599
$ws.Range("I106").Value = 6498.75
$ws.Range("K106").Value = 6498.75
$ws.Range("M106").Value = -5867.75

$ws.Range("H132").Value = 5450.25
$ws.Range("I132").Value = 4059.2727
$ws.Range("K132").Value = 12177.8181
$ws.Range("M132").Value = -9647.8181

$ws.Range("H138").Value = 3880.92
$ws.Range("I138").Value = 3785.4
$ws.Range("J138").Value = 3944.6
$ws.Range("K138").Value = 11356.2
$ws.Range("L138").Value = 11833.8
$ws.Range("M138").Value = -6216.200000000001
$ws.Range("N138").Value = -22113.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3418.8
$ws.Range("I61").Value = 3418.8
$ws.Range("K61").Value = 3418.8
$ws.Range("M61").Value = -3206.8

$ws.Range("H74").Value = 4166.6665
$ws.Range("I74").Value = 5250
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 5250
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -4376
$ws.Range("N74").Value = -3748

$ws.Range("H77").Value = 4166.6665
$ws.Range("I77").Value = 5250
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 26250
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -21882
$ws.Range("N77").Value = -18736

$ws.Range("H136").Value = 3418.8
$ws.Range("I136").Value = 3418.8
$ws.Range("K136").Value = 10256.4
$ws.Range("M136").Value = -7706.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4250
$ws.Range("I20").Value = 2500
$ws.Range("J20").Value = 6000
$ws.Range("K20").Value = 2500
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = -2253
$ws.Range("N20").Value = -6494

$ws.Range("H22").Value = 999.75
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 999
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -826
$ws.Range("N22").Value = -1346

$ws.Range("H99").Value = 4099.8
$ws.Range("I99").Value = 4099.8
$ws.Range("K99").Value = 4099.8
$ws.Range("M99").Value = -2601.8

$ws.Range("H134").Value = 6994.3335
$ws.Range("I134").Value = 4278.7144
$ws.Range("K134").Value = 12836.1432
$ws.Range("M134").Value = -10301.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23497.875
$ws.Range("I31").Value = 15097.5
$ws.Range("J31").Value = 37498.5
$ws.Range("K31").Value = 15097.5
$ws.Range("L31").Value = 37498.5
$ws.Range("M31").Value = -14802.5
$ws.Range("N31").Value = -38088.5

$ws.Range("H34").Value = 23497.875
$ws.Range("I34").Value = 15097.5
$ws.Range("J34").Value = 37498.5
$ws.Range("K34").Value = 15097.5
$ws.Range("L34").Value = 37498.5
$ws.Range("M34").Value = -14895.5
$ws.Range("N34").Value = -37902.5

$ws.Range("H122").Value = 2979
$ws.Range("I122").Value = 2979
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8937
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -6487

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 2407.3333
$ws.Range("I23").Value = 2111
$ws.Range("J23").Value = 3000
$ws.Range("K23").Value = 6333
$ws.Range("L23").Value = 9000
$ws.Range("M23").Value = -6098
$ws.Range("N23").Value = -9470

$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()

$ws.Range("H113").Value = 995.6667
$ws.Range("I113").Value = 737.5
$ws.Range("K113").Value = 2212.5
$ws.Range("M113").Value = -42.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1431971.2
$ws.Range("I14").Value = 3949.5
$ws.Range("J14").Value = 2003180
$ws.Range("K14").Value = 3949.5
$ws.Range("L14").Value = 2003180
$ws.Range("M14").Value = -3781.5
$ws.Range("N14").Value = -2003516

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 841.4286
$ws.Range("I22").Value = 725
$ws.Range("K22").Value = 725
$ws.Range("M22").Value = -430

$ws.Range("H27").Value = 841.4286
$ws.Range("I27").Value = 725
$ws.Range("K27").Value = 725
$ws.Range("M27").Value = -618

$ws.Range("H46").Value = 867.8570999999999
$ws.Range("I46").Value = 775
$ws.Range("J46").Value = 937.5
$ws.Range("K46").Value = 775
$ws.Range("L46").Value = 937.5
$ws.Range("M46").Value = -587
$ws.Range("N46").Value = -1313.5

$ws.Range("H61").Value = 6913.6665
$ws.Range("I61").Value = 6898.4
$ws.Range("J61").Value = 6990
$ws.Range("K61").Value = 6898.4
$ws.Range("L61").Value = 6990
$ws.Range("M61").Value = -6696.4
$ws.Range("N61").Value = -7394

$ws.Range("H68").Value = 2880.3635
$ws.Range("I68").Value = 2824.625
$ws.Range("J68").Value = 3029
$ws.Range("K68").Value = 2824.625
$ws.Range("L68").Value = 3029
$ws.Range("M68").Value = -2075.625
$ws.Range("N68").Value = -4527

$ws.Range("H71").Value = 2880.3635
$ws.Range("I71").Value = 2824.625
$ws.Range("J71").Value = 3029
$ws.Range("K71").Value = 14123.125
$ws.Range("L71").Value = 15145
$ws.Range("M71").Value = -10379.125
$ws.Range("N71").Value = -22633

$ws.Range("H74").Value = 50000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 50000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H113").Value = 6913.6665
$ws.Range("I113").Value = 6898.4
$ws.Range("J113").Value = 6990
$ws.Range("K113").Value = 6898.4
$ws.Range("L113").Value = 6990
$ws.Range("M113").Value = -4728.4
$ws.Range("N113").Value = -11330

$ws.Range("H132").Value = 29400
$ws.Range("I132").Value = 41750
$ws.Range("J132").Value = 17050
$ws.Range("K132").Value = 125250
$ws.Range("L132").Value = 51150
$ws.Range("M132").Value = -122720
$ws.Range("N132").Value = -56210

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0

$ws.Range("H126").Value = 2270.8572
$ws.Range("I126").Value = 1149.25
$ws.Range("K126").Value = 3447.75
$ws.Range("M126").Value = -977.75

$ws.Range("H132").Value = 9498.5
$ws.Range("I132").Value = 5998.5
$ws.Range("K132").Value = 17995.5
$ws.Range("M132").Value = -15465.5
